$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7602039999999999
$ws.Range("H2").Value = 2.280612
$ws.Range("I2").Value = 0.9081302434927829
$ws.Range("J2").Value = 0.9140443705567521
$ws.Range("M2").Value = 0.6598136666666666
$ws.Range("N2").Value = 1.979441
$ws.Range("O2").Value = 0.007704735356083927
$ws.Range("P2").Value = 0.008484678519943686
$ws.Range("Q2").Value = 0.5015929886546666
$ws.Range("R2").Value = 4.514336897891999
$ws.Range("S2").Value = 0.00699690319496795
$ws.Range("T2").Value = 0.007755372637138322

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7602039999999999
$ws.Range("H3").Value = 2.280612
$ws.Range("I3").Value = 0.9081302434927829
$ws.Range("J3").Value = 0.9140443705567521
$ws.Range("O3").Value = 0.7130079175842846
$ws.Range("P3").Value = 0.7851850431306702
$ws.Range("Q3").Value = 46.41817736583732
$ws.Range("R3").Value = 417.7635962925359
$ws.Range("S3").Value = 0.6475040538080985
$ws.Range("T3").Value = 0.7176939685189497

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7602039999999999
$ws.Range("H4").Value = 2.280612
$ws.Range("I4").Value = 0.9081302434927829
$ws.Range("J4").Value = 0.9140443705567521
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1563486666666667
$ws.Range("N4").Value = 0.469046
$ws.Range("O4").Value = 0.001825704984300993
$ws.Range("P4").Value = 0.002010519394650058
$ws.Range("Q4").Value = 0.1188568817946666
$ws.Range("R4").Value = 1.069711936152
$ws.Range("S4").Value = 0.001657977911939248
$ws.Range("T4").Value = 0.001837703934575055

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7602039999999999
$ws.Range("H5").Value = 2.280612
$ws.Range("I5").Value = 0.9081302434927829
$ws.Range("J5").Value = 0.9140443705567521
$ws.Range("M5").Value = 23.6163295
$ws.Range("N5").Value = 47.232659
$ws.Range("O5").Value = 0.2757711427815902
$ws.Range("P5").Value = 0.2024581319964196
$ws.Range("Q5").Value = 17.953228151218
$ws.Range("R5").Value = 107.719368907308
$ws.Range("S5").Value = 0.2504361150425285
$ws.Range("T5").Value = 0.1850557158247632

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7602039999999999
$ws.Range("H6").Value = 2.280612
$ws.Range("I6").Value = 0.9081302434927829
$ws.Range("J6").Value = 0.9140443705567521
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.14477
$ws.Range("N6").Value = 0.43431
$ws.Range("O6").Value = 0.00169049929374041
$ws.Range("P6").Value = 0.001861626958316384
$ws.Range("Q6").Value = 0.11005473308
$ws.Range("R6").Value = 0.9904925977199999
$ws.Range("S6").Value = 0.001535193535248856
$ws.Range("T6").Value = 0.001701609641325781

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.016249
$ws.Range("H7").Value = 0.032498
$ws.Range("I7").Value = 0.01941085330584189
$ws.Range("J7").Value = 0.01302484331151171
$ws.Range("M7").Value = 0.6598136666666666
$ws.Range("N7").Value = 1.979441
$ws.Range("O7").Value = 0.007704735356083927
$ws.Range("P7").Value = 0.008484678519943686
$ws.Range("Q7").Value = 0.01072131226966667
$ws.Range("R7").Value = 0.06432787361799999
$ws.Range("S7").Value = 0.0001495554877572786
$ws.Range("T7").Value = 0.0001105116082708156

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.016249
$ws.Range("H8").Value = 0.032498
$ws.Range("I8").Value = 0.01941085330584189
$ws.Range("J8").Value = 0.01302484331151171
$ws.Range("O8").Value = 0.7130079175842846
$ws.Range("P8").Value = 0.7851850431306702
$ws.Range("Q8").Value = 0.9921665290073333
$ws.Range("R8").Value = 5.952999174044
$ws.Range("S8").Value = 0.01384009209413235
$ws.Range("T8").Value = 0.01022691215731954

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.016249
$ws.Range("H9").Value = 0.032498
$ws.Range("I9").Value = 0.01941085330584189
$ws.Range("J9").Value = 0.01302484331151171
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1563486666666667
$ws.Range("N9").Value = 0.469046
$ws.Range("O9").Value = 0.001825704984300993
$ws.Range("P9").Value = 0.002010519394650058
$ws.Range("Q9").Value = 0.002540509484666666
$ws.Range("R9").Value = 0.015243056908
$ws.Range("S9").Value = 0.00003543849163001095
$ws.Range("T9").Value = 0.00002618670009007238

# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.016249
$ws.Range("H10").Value = 0.032498
$ws.Range("I10").Value = 0.01941085330584189
$ws.Range("J10").Value = 0.01302484331151171
$ws.Range("M10").Value = 23.6163295
$ws.Range("N10").Value = 47.232659
$ws.Range("O10").Value = 0.2757711427815902
$ws.Range("P10").Value = 0.2024581319964196
$ws.Range("Q10").Value = 0.3837417380455
$ws.Range("R10").Value = 1.534966952182
$ws.Range("S10").Value = 0.005352953198517826
$ws.Range("T10").Value = 0.00263698544639472

# Row 11
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.016249
$ws.Range("H11").Value = 0.032498
$ws.Range("I11").Value = 0.01941085330584189
$ws.Range("J11").Value = 0.01302484331151171
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.14477
$ws.Range("N11").Value = 0.43431
$ws.Range("O11").Value = 0.00169049929374041
$ws.Range("P11").Value = 0.001861626958316384
$ws.Range("Q11").Value = 0.00235236773
$ws.Range("R11").Value = 0.01411420638
$ws.Range("S11").Value = 0.00003281403380442442
$ws.Range("T11").Value = 0.00002424739943655704

# Row 12
$ws.Range("G12").Value = 0.06065600000000001
$ws.Range("H12").Value = 0.181968
$ws.Range("I12").Value = 0.07245890320137523
$ws.Range("J12").Value = 0.07293078613173619
$ws.Range("M12").Value = 0.6598136666666666
$ws.Range("N12").Value = 1.979441
$ws.Range("O12").Value = 0.007704735356083927
$ws.Range("P12").Value = 0.008484678519943686
$ws.Range("Q12").Value = 0.04002165776533333
$ws.Range("R12").Value = 0.3601949198880001
$ws.Range("S12").Value = 0.0005582766733586986
$ws.Range("T12").Value = 0.0006187942745345489

# Row 13
$ws.Range("G13").Value = 0.06065600000000001
$ws.Range("H13").Value = 0.181968
$ws.Range("I13").Value = 0.07245890320137523
$ws.Range("J13").Value = 0.07293078613173619
$ws.Range("O13").Value = 0.7130079175842846
$ws.Range("P13").Value = 0.7851850431306702
$ws.Range("Q13").Value = 3.703665024522667
$ws.Range("R13").Value = 33.332985220704
$ws.Range("S13").Value = 0.05166377168205381
$ws.Range("T13").Value = 0.05726416245440096

# Row 14
$ws.Range("G14").Value = 0.06065600000000001
$ws.Range("H14").Value = 0.181968
$ws.Range("I14").Value = 0.07245890320137523
$ws.Range("J14").Value = 0.07293078613173619
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1563486666666667
$ws.Range("N14").Value = 0.469046
$ws.Range("O14").Value = 0.001825704984300993
$ws.Range("P14").Value = 0.002010519394650058
$ws.Range("Q14").Value = 0.009483484725333334
$ws.Range("R14").Value = 0.085351362528
$ws.Range("S14").Value = 0.0001322885807317339
$ws.Range("T14").Value = 0.0001466287599849311

# Row 15
$ws.Range("G15").Value = 0.06065600000000001
$ws.Range("H15").Value = 0.181968
$ws.Range("I15").Value = 0.07245890320137523
$ws.Range("J15").Value = 0.07293078613173619
$ws.Range("M15").Value = 23.6163295
$ws.Range("N15").Value = 47.232659
$ws.Range("O15").Value = 0.2757711427815902
$ws.Range("P15").Value = 0.2024581319964196
$ws.Range("Q15").Value = 1.432472082152
$ws.Range("R15").Value = 8.594832492912001
$ws.Range("S15").Value = 0.01998207454054387
$ws.Range("T15").Value = 0.01476543072526169

# Row 16
$ws.Range("G16").Value = 0.06065600000000001
$ws.Range("H16").Value = 0.181968
$ws.Range("I16").Value = 0.07245890320137523
$ws.Range("J16").Value = 0.07293078613173619
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.14477
$ws.Range("N16").Value = 0.43431
$ws.Range("O16").Value = 0.00169049929374041
$ws.Range("P16").Value = 0.001861626958316384
$ws.Range("Q16").Value = 0.008781169120000001
$ws.Range("R16").Value = 0.07903052208000001
$ws.Range("S16").Value = 0.0001224917246871295
$ws.Range("T16").Value = 0.0001357699175540468
